# Auto-generated edits applying the Aegis_Profits.xlsx market-data refresh
# (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 5797.143
$ws.Range("I40").Value = 7120
$ws.Range("J40").Value = 2490
$ws.Range("K40").Value = 7120
$ws.Range("L40").Value = 2490
$ws.Range("M40").Value = -6945
$ws.Range("N40").Value = -2840

# Row 64
$ws.Range("H64").Value = 65502.688
$ws.Range("I64").Value = 169110
$ws.Range("J64").Value = 3338.3
$ws.Range("K64").Value = 169110
$ws.Range("L64").Value = 3338.3
$ws.Range("M64").Value = -168862
$ws.Range("N64").Value = -3834.3

# Row 67
$ws.Range("H67").Value = 65502.688
$ws.Range("I67").Value = 169110
$ws.Range("J67").Value = 3338.3
$ws.Range("K67").Value = 169110
$ws.Range("L67").Value = 3338.3
$ws.Range("M67").Value = -168252
$ws.Range("N67").Value = -5054.3

# Row 129
$ws.Range("H129").Value = 2316.4707
$ws.Range("I129").Value = 6753.625
$ws.Range("J129").Value = 951.1923
$ws.Range("K129").Value = 20260.875
$ws.Range("L129").Value = 2853.5769
$ws.Range("M129").Value = -15260.875
$ws.Range("N129").Value = -12853.5769

# Row 138
$ws.Range("H138").Value = 3452.3394
$ws.Range("I138").Value = 7333.3335
$ws.Range("J138").Value = 3232.6604
$ws.Range("K138").Value = 22000.0005
$ws.Range("L138").Value = 9697.9812
$ws.Range("M138").Value = -16860.0005
$ws.Range("N138").Value = -19977.9812

# Row 141
$ws.Range("H141").Value = 2041.2106
$ws.Range("I141").Value = 1478.5333
$ws.Range("K141").Value = 4435.5999
$ws.Range("M141").Value = 744.4000999999998

$ws = $wb.Worksheets.Item("ARM")
# Row 25
$ws.Range("H25").Value = 5600
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 5600
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 5600
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -6404

# Row 74
$ws.Range("H74").Value = 1021.1905
$ws.Range("I74").Value = 947.3570999999999
$ws.Range("K74").Value = 947.3570999999999
$ws.Range("M74").Value = -73.35709999999995

# Row 77
$ws.Range("H77").Value = 1021.1905
$ws.Range("I77").Value = 947.3570999999999
$ws.Range("K77").Value = 4736.7855
$ws.Range("M77").Value = -368.7855

# Row 132
$ws.Range("H132").Value = 12625.075
$ws.Range("I132").Value = 15434.317
$ws.Range("J132").Value = 3026.8333
$ws.Range("K132").Value = 46302.951
$ws.Range("L132").Value = 9080.499899999999
$ws.Range("M132").Value = -43772.951
$ws.Range("N132").Value = -14140.4999

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 13073.341
$ws.Range("I134").Value = 14633.7
$ws.Range("J134").Value = 4157
$ws.Range("K134").Value = 43901.10000000001
$ws.Range("L134").Value = 12471
$ws.Range("M134").Value = -41366.10000000001
$ws.Range("N134").Value = -17541

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 36711.242
$ws.Range("I31").Value = 643.4211
$ws.Range("K31").Value = 643.4211
$ws.Range("M31").Value = -348.4211

# Row 34
$ws.Range("H34").Value = 36711.242
$ws.Range("I34").Value = 643.4211
$ws.Range("K34").Value = 643.4211
$ws.Range("M34").Value = -441.4211

# Row 58
$ws.Range("H58").Value = 2481.077
$ws.Range("I58").Value = 2214
$ws.Range("J58").Value = 3371.3333
$ws.Range("K58").Value = 2214
$ws.Range("L58").Value = 3371.3333
$ws.Range("M58").Value = -2011
$ws.Range("N58").Value = -3777.3333

# Row 136
$ws.Range("H136").Value = 2481.077
$ws.Range("I136").Value = 2214
$ws.Range("J136").Value = 3371.3333
$ws.Range("K136").Value = 6642
$ws.Range("L136").Value = 10113.9999
$ws.Range("M136").Value = -4092
$ws.Range("N136").Value = -15213.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1250.4833
$ws.Range("I5").Value = 1272.0834
$ws.Range("J5").Value = 1245.0834
$ws.Range("K5").Value = 3816.2502
$ws.Range("L5").Value = 3735.2502
$ws.Range("M5").Value = -3704.2502
$ws.Range("N5").Value = -3959.2502

# Row 12
$ws.Range("H12").Value = 43.588234
$ws.Range("I12").Value = 15.666667
$ws.Range("J12").Value = 49.57143
$ws.Range("K12").Value = 47.000001
$ws.Range("L12").Value = 148.71429
$ws.Range("M12").Value = 125.999999
$ws.Range("N12").Value = -494.71429

# Row 47
$ws.Range("H47").Value = 200
$ws.Range("I47").Value = 200
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 600
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -169
$ws.Range("N47").ClearContents()

# Row 107
$ws.Range("H107").Value = 935.6774
$ws.Range("I107").Value = 989.6
$ws.Range("J107").Value = 925.3077
$ws.Range("K107").Value = 2968.8
$ws.Range("L107").Value = 2775.9231
$ws.Range("M107").Value = -1048.8
$ws.Range("N107").Value = -6615.9231

# Row 131
$ws.Range("H131").Value = 17515.895
$ws.Range("I131").Value = 479.9091
$ws.Range("J131").Value = 20923.092
$ws.Range("K131").Value = 1439.7273
$ws.Range("L131").Value = 62769.276
$ws.Range("M131").Value = 3600.2727
$ws.Range("N131").Value = -72849.276

# Row 132
$ws.Range("H132").Value = 402033.9
$ws.Range("I132").Value = 553.36365
$ws.Range("J132").Value = 717482.9399999999
$ws.Range("K132").Value = 4980.27285
$ws.Range("L132").Value = 6457346.459999999
$ws.Range("M132").Value = -2450.27285
$ws.Range("N132").Value = -6462406.459999999

# Row 135
$ws.Range("H135").Value = 1250.4833
$ws.Range("I135").Value = 1272.0834
$ws.Range("J135").Value = 1245.0834
$ws.Range("K135").Value = 11448.7506
$ws.Range("L135").Value = 11205.7506
$ws.Range("M135").Value = -8913.750599999999
$ws.Range("N135").Value = -16275.7506

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2643.3333
$ws.Range("I7").Value = 1686.4
$ws.Range("J7").Value = 3839.5
$ws.Range("K7").Value = 1686.4
$ws.Range("L7").Value = 3839.5
$ws.Range("M7").Value = -1574.4
$ws.Range("N7").Value = -4063.5

# Row 93
$ws.Range("H93").Value = 2224.36
$ws.Range("I93").Value = 2163.611
$ws.Range("J93").Value = 2380.5715
$ws.Range("K93").Value = 2163.611
$ws.Range("L93").Value = 2380.5715
$ws.Range("M93").Value = -915.6109999999999
$ws.Range("N93").Value = -4876.5715

# Row 126
$ws.Range("H126").Value = 2643.3333
$ws.Range("I126").Value = 1686.4
$ws.Range("J126").Value = 3839.5
$ws.Range("K126").Value = 5059.200000000001
$ws.Range("L126").Value = 11518.5
$ws.Range("M126").Value = -2589.200000000001
$ws.Range("N126").Value = -16458.5

# Row 131
$ws.Range("H131").Value = 29000
$ws.Range("J131").Value = 29000
$ws.Range("L131").Value = 29000
$ws.Range("N131").Value = -39080

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 10991256
$ws.Range("I62").Value = 76923070
$ws.Range("J62").Value = 2620.3333
$ws.Range("K62").Value = 76923070
$ws.Range("L62").Value = 2620.3333
$ws.Range("M62").Value = -76922446
$ws.Range("N62").Value = -3868.3333

# Row 65
$ws.Range("H65").Value = 10991256
$ws.Range("I65").Value = 76923070
$ws.Range("J65").Value = 2620.3333
$ws.Range("K65").Value = 384615350
$ws.Range("L65").Value = 13101.6665
$ws.Range("M65").Value = -384612230
$ws.Range("N65").Value = -19341.6665

# Row 76
$ws.Range("H76").Value = 33000
$ws.Range("J76").Value = 33000
$ws.Range("L76").Value = 33000
$ws.Range("N76").Value = -33630

# Row 79
$ws.Range("H79").Value = 33000
$ws.Range("J79").Value = 33000
$ws.Range("L79").Value = 33000
$ws.Range("N79").Value = -35184

# Row 122
$ws.Range("H122").Value = 2305.0435
$ws.Range("I122").Value = 1882.7059
$ws.Range("J122").Value = 3501.6667
$ws.Range("K122").Value = 5648.1177
$ws.Range("L122").Value = 10505.0001
$ws.Range("M122").Value = -3198.1177
$ws.Range("N122").Value = -15405.0001

# Row 132
$ws.Range("H132").Value = 4249.755
$ws.Range("I132").Value = 2630.8965
$ws.Range("J132").Value = 6597.1
$ws.Range("K132").Value = 7892.689499999999
$ws.Range("L132").Value = 19791.3
$ws.Range("M132").Value = -5362.689499999999
$ws.Range("N132").Value = -24851.3

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
